$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnc"
$ws.Range("C2").Value = "Itgb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 8.948174
$ws.Range("H2").Value = 26.844522
$ws.Range("I2").Value = 0.0695931738232498
$ws.Range("J2").Value = 0.0695931738232498
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 3.778439
$ws.Range("N2").Value = 11.335317
$ws.Range("O2").Value = 0.4252971528324392
$ws.Range("P2").Value = 0.4252971528324392
$ws.Range("Q2").Value = 33.810129620386
$ws.Range("R2").Value = 304.291166583474
$ws.Range("S2").Value = 0.02959777868360118
$ws.Range("T2").Value = 0.02959777868360118

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tnc"
$ws.Range("C3").Value = "Itgb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 8.948174
$ws.Range("H3").Value = 26.844522
$ws.Range("I3").Value = 0.0695931738232498
$ws.Range("J3").Value = 0.0695931738232498
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("N3").Value = 13.00021
$ws.Range("O3").Value = 0.4877633593505858
$ws.Range("P3").Value = 0.4877633593505858
$ws.Range("Q3").Value = 38.77604703884666
$ws.Range("R3").Value = 348.9844233496199
$ws.Range("S3").Value = 0.03394500025189757
$ws.Range("T3").Value = 0.03394500025189758

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tnc"
$ws.Range("C4").Value = "Itgb3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 8.948174
$ws.Range("H4").Value = 26.844522
$ws.Range("I4").Value = 0.0695931738232498
$ws.Range("J4").Value = 0.0695931738232498
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 0.2909853333333334
$ws.Range("N4").Value = 0.8729560000000001
$ws.Range("O4").Value = 0.03275300561492853
$ws.Range("P4").Value = 0.03275300561492853
$ws.Range("Q4").Value = 2.603787394114667
$ws.Range("R4").Value = 23.434086547032
$ws.Range("S4").Value = 0.002279385612993598
$ws.Range("T4").Value = 0.002279385612993598

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Tnc"
$ws.Range("C5").Value = "Itgb3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 8.948174
$ws.Range("H5").Value = 26.844522
$ws.Range("I5").Value = 0.0695931738232498
$ws.Range("J5").Value = 0.0695931738232498
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 0.4814053333333333
$ws.Range("N5").Value = 1.444216
$ws.Range("O5").Value = 0.0541864822020464
$ws.Range("P5").Value = 0.05418648220204641
$ws.Range("Q5").Value = 4.307698687194666
$ws.Range("R5").Value = 38.76928818475199
$ws.Range("S5").Value = 0.003771009274757446
$ws.Range("T5").Value = 0.003771009274757447

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tnc"
$ws.Range("C6").Value = "Itgb3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 105.8801323333333
$ws.Range("H6").Value = 317.640397
$ws.Range("I6").Value = 0.8234679448457706
$ws.Range("J6").Value = 0.8234679448457706
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 3.778439
$ws.Range("N6").Value = 11.335317
$ws.Range("O6").Value = 0.4252971528324392
$ws.Range("P6").Value = 0.4252971528324392
$ws.Range("Q6").Value = 400.0616213334277
$ws.Range("R6").Value = 3600.554592000849
$ws.Range("S6").Value = 0.3502185723916864
$ws.Range("T6").Value = 0.3502185723916864

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tnc"
$ws.Range("C7").Value = "Itgb3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 105.8801323333333
$ws.Range("H7").Value = 317.640397
$ws.Range("I7").Value = 0.8234679448457706
$ws.Range("J7").Value = 0.8234679448457706
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("N7").Value = 13.00021
$ws.Range("O7").Value = 0.4877633593505858
$ws.Range("P7").Value = 0.4877633593505858
$ws.Range("Q7").Value = 458.8213183870411
$ws.Range("R7").Value = 4129.39186548337
$ws.Range("S7").Value = 0.401657491095496
$ws.Range("T7").Value = 0.401657491095496

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Tnc"
$ws.Range("C8").Value = "Itgb3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 105.8801323333333
$ws.Range("H8").Value = 317.640397
$ws.Range("I8").Value = 0.8234679448457706
$ws.Range("J8").Value = 0.8234679448457706
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 0.2909853333333334
$ws.Range("N8").Value = 0.8729560000000001
$ws.Range("O8").Value = 0.03275300561492853
$ws.Range("P8").Value = 0.03275300561492853
$ws.Range("Q8").Value = 30.80956560039245
$ws.Range("R8").Value = 277.286090403532
$ws.Range("S8").Value = 0.02697105022124718
$ws.Range("T8").Value = 0.02697105022124718

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Tnc"
$ws.Range("C9").Value = "Itgb3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 105.8801323333333
$ws.Range("H9").Value = 317.640397
$ws.Range("I9").Value = 0.8234679448457706
$ws.Range("J9").Value = 0.8234679448457706
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 0.4814053333333333
$ws.Range("N9").Value = 1.444216
$ws.Range("O9").Value = 0.0541864822020464
$ws.Range("P9").Value = 0.05418648220204641
$ws.Range("Q9").Value = 50.97126039930578
$ws.Range("R9").Value = 458.741343593752
$ws.Range("S9").Value = 0.04462083113734108
$ws.Range("T9").Value = 0.04462083113734108

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Tnc"
$ws.Range("C10").Value = "Itgb3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2.0
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1148696666666667
$ws.Range("H10").Value = 0.344609
$ws.Range("I10").Value = 0.0008933827928862465
$ws.Range("J10").Value = 0.0008933827928862465
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 3.778439
$ws.Range("N10").Value = 11.335317
$ws.Range("O10").Value = 0.4252971528324392
$ws.Range("P10").Value = 0.4252971528324392
$ws.Range("Q10").Value = 0.4340280284503333
$ws.Range("R10").Value = 3.906252256053
$ws.Range("S10").Value = 0.0003799531582040134
$ws.Range("T10").Value = 0.0003799531582040134

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Tnc"
$ws.Range("C11").Value = "Itgb3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2.0
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1148696666666667
$ws.Range("H11").Value = 0.344609
$ws.Range("I11").Value = 0.0008933827928862465
$ws.Range("J11").Value = 0.0008933827928862465
$ws.Range("K11").Value = 3.0
$ws.Range("L11").Value = 1.0
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("N11").Value = 13.00021
$ws.Range("O11").Value = 0.4877633593505858
$ws.Range("P11").Value = 0.4877633593505858
$ws.Range("Q11").Value = 0.4977765964322222
$ws.Range("R11").Value = 4.47998936789
$ws.Range("S11").Value = 0.0004357593922442042
$ws.Range("T11").Value = 0.0004357593922442043

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Tnc"
$ws.Range("C12").Value = "Itgb3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 2.0
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1148696666666667
$ws.Range("H12").Value = 0.344609
$ws.Range("I12").Value = 0.0008933827928862465
$ws.Range("J12").Value = 0.0008933827928862465
$ws.Range("K12").Value = 3.0
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 0.2909853333333334
$ws.Range("N12").Value = 0.8729560000000001
$ws.Range("O12").Value = 0.03275300561492853
$ws.Range("P12").Value = 0.03275300561492853
$ws.Range("Q12").Value = 0.03342538824488889
$ws.Range("R12").Value = 0.300828494204
$ws.Range("S12").Value = 0.00002926097163168377
$ws.Range("T12").Value = 0.00002926097163168377

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Tnc"
$ws.Range("C13").Value = "Itgb3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 2.0
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1148696666666667
$ws.Range("H13").Value = 0.344609
$ws.Range("I13").Value = 0.0008933827928862465
$ws.Range("J13").Value = 0.0008933827928862465
$ws.Range("K13").Value = 3.0
$ws.Range("L13").Value = 1.0
$ws.Range("M13").Value = 0.4814053333333333
$ws.Range("N13").Value = 1.444216
$ws.Range("O13").Value = 0.0541864822020464
$ws.Range("P13").Value = 0.05418648220204641
$ws.Range("Q13").Value = 0.05529887017155555
$ws.Range("R13").Value = 0.497689831544
$ws.Range("S13").Value = 0.00004840927080634511
$ws.Range("T13").Value = 0.00004840927080634511

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Tnc"
$ws.Range("C14").Value = "Itgb3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 13.635153
$ws.Range("H14").Value = 40.905459
$ws.Range("I14").Value = 0.1060454985380935
$ws.Range("J14").Value = 0.1060454985380935
$ws.Range("K14").Value = 3.0
$ws.Range("L14").Value = 1.0
$ws.Range("M14").Value = 3.778439
$ws.Range("N14").Value = 11.335317
$ws.Range("O14").Value = 0.4252971528324392
$ws.Range("P14").Value = 0.4252971528324392
$ws.Range("Q14").Value = 51.51959386616701
$ws.Range("R14").Value = 463.676344795503
$ws.Range("S14").Value = 0.04510084859894775
$ws.Range("T14").Value = 0.04510084859894774

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Tnc"
$ws.Range("C15").Value = "Itgb3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 13.635153
$ws.Range("H15").Value = 40.905459
$ws.Range("I15").Value = 0.1060454985380935
$ws.Range("J15").Value = 0.1060454985380935
$ws.Range("K15").Value = 3.0
$ws.Range("L15").Value = 1.0
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("N15").Value = 13.00021
$ws.Range("O15").Value = 0.4877633593505858
$ws.Range("P15").Value = 0.4877633593505858
$ws.Range("Q15").Value = 59.08661746071
$ws.Range("R15").Value = 531.77955714639
$ws.Range("S15").Value = 0.05172510861094811
$ws.Range("T15").Value = 0.0517251086109481

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Tnc"
$ws.Range("C16").Value = "Itgb3"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 13.635153
$ws.Range("H16").Value = 40.905459
$ws.Range("I16").Value = 0.1060454985380935
$ws.Range("J16").Value = 0.1060454985380935
$ws.Range("K16").Value = 3.0
$ws.Range("L16").Value = 1.0
$ws.Range("M16").Value = 0.2909853333333334
$ws.Range("N16").Value = 0.8729560000000001
$ws.Range("O16").Value = 0.03275300561492853
$ws.Range("P16").Value = 0.03275300561492853
$ws.Range("Q16").Value = 3.967629540756001
$ws.Range("R16").Value = 35.708665866804
$ws.Range("S16").Value = 0.003473308809056071
$ws.Range("T16").Value = 0.003473308809056071

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Tnc"
$ws.Range("C17").Value = "Itgb3"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3.0
$ws.Range("F17").Value = 1.0
$ws.Range("G17").Value = 13.635153
$ws.Range("H17").Value = 40.905459
$ws.Range("I17").Value = 0.1060454985380935
$ws.Range("J17").Value = 0.1060454985380935
$ws.Range("K17").Value = 3.0
$ws.Range("L17").Value = 1.0
$ws.Range("M17").Value = 0.4814053333333333
$ws.Range("N17").Value = 1.444216
$ws.Range("O17").Value = 0.0541864822020464
$ws.Range("P17").Value = 0.05418648220204641
$ws.Range("Q17").Value = 6.564035375016
$ws.Range("R17").Value = 59.076318375144
$ws.Range("S17").Value = 0.00574623251914154
$ws.Range("T17").Value = 0.00574623251914154
